# "Generate Report for Handback"
#
# The handback transform for 839376d6-c616-480e-963b-bc5917dce785.md failed
# (handback type "mt" didn't match the handoff type "ht"), so its row needs
# to move from "Ready for handoff" to "Handback transform failed" with an
# Error Detail message, for both locale sheets. The Overview sheet's
# per-locale status columns for that same file are updated to match.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$errorDetail = "The handback type mt is not match with handoff type ht."

# Overview sheet: row for 839376d6-c616-480e-963b-bc5917dce785.md is row 6;
# column B = zh-cn status, column C = de-de status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B6").Value = $newStatus
$wsOverview.Range("C6").Value = $newStatus

# zh-cn detail sheet: same file is row 6; column C = Status, column L = Error Detail.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = $newStatus
$wsZhCn.Range("L6").Value = $errorDetail

# de-de detail sheet: same file is row 6; column C = Status, column L = Error Detail.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = $newStatus
$wsDeDe.Range("L6").Value = $errorDetail
